$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Helper: apply the "h:mm" time style (matches existing style index 5) to a cell
function Set-TimeCell($row, $col, $value) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "h:mm"
    $c.Value = $value
}

# --- Row 20: 10-3-2017 (date 42804), Software Architect, "Lavet design af OC3" ---
$ws.Cells.Item(20, 1).Value = 42804
$ws.Cells.Item(20, 5).Value = "Software Architect"
$ws.Cells.Item(20, 6).Value = "Lavet design af OC3"
Set-TimeCell 20 7 0.33680555555555558
Set-TimeCell 20 8 0.43055555555555558

# --- Row 21: Requirements Specifier, "Begynd at lave test suite for OC4 " ---
$ws.Cells.Item(21, 5).Value = "Requirements Specifier"
$ws.Cells.Item(21, 6).Value = "Begynd at lave test suite for OC4 "
Set-TimeCell 21 7 0.5
Set-TimeCell 21 8 0.52083333333333337

# --- Row 22: Reviewer, "Lavet review af OC artefakter" ---
$ws.Cells.Item(22, 5).Value = "Reviewer"
$ws.Cells.Item(22, 6).Value = "Lavet review af OC artefakter"
Set-TimeCell 22 7 0.43055555555555558
Set-TimeCell 22 8 0.47916666666666669

# --- Row 23: total hours for the day ---
$ws.Cells.Item(23, 9).Value = 3.55

# --- Row 24: 13-3-2017 (date 42807), Reviewer, "Lavet review af test suite for OC4 " ---
$ws.Cells.Item(24, 1).Value = 42807
# Reuse the date style (s="1") from an existing date-formatted column A cell
$ws.Cells.Item(21, 1).Copy()
$ws.Cells.Item(24, 1).PasteSpecial(-4122)
$ws.Cells.Item(24, 5).Value = "Reviewer"
$ws.Cells.Item(24, 6).Value = "Lavet review af test suite for OC4 "
Set-TimeCell 24 7 0.34861111111111115
Set-TimeCell 24 8 0.38541666666666669

# --- Row 26: "Lavet review af OC4" (shared-string text entered before row 25's) ---
$ws.Cells.Item(26, 6).Value = "Lavet review af OC4"

# --- Row 25: "Lavet review af OC2" ---
$ws.Cells.Item(25, 6).Value = "Lavet review af OC2"
Set-TimeCell 25 7 0.4375
Set-TimeCell 25 8 0.47222222222222227

Set-TimeCell 26 7 0.50486111111111109
Set-TimeCell 26 8 0.56527777777777777

# --- Row 27: Requirements Specifier, "Integreret Dom5 " ---
$ws.Cells.Item(27, 5).Value = "Requirements Specifier"
$ws.Cells.Item(27, 6).Value = "Integreret Dom5 "
Set-TimeCell 27 7 0.61458333333333337
Set-TimeCell 27 8 0.62222222222222223

# --- Row 28: total hours for the day ---
$ws.Cells.Item(28, 9).Value = 3.2

# --- Update the view: scroll down and select E29 ---
$ws.Range("E29").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 3
$win.ScrollColumn = 1
